$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new centered heading paragraph "Báo cáo hằng tuần" right
#    before the existing "Đồ án ... " title paragraph (paragraph 5), reusing
#    that paragraph's formatting (InsertParagraphBefore clones pPr/rPr).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(5)
$titlePara.Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs(5).Range.Text = "Báo cáo hằng tuần"

# ---------------------------------------------------------------------------
# 2) Remove one of the two consecutive empty paragraphs that follow
#    "Hệ thống lưu thói quen tắt đèn của người dùng để tự động tắt." —
#    keep only a single blank paragraph afterwards.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "Hệ thống lưu thói quen tắt đèn của người dùng để tự động tắt\.") {
        $d.Paragraphs($i + 1).Range.Delete() | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Fix the typo "tiềm hiểu" -> "tìm hiểu" in both occurrences of
#    "Đã tiềm hiểu sơ bộ, đang tiến hành thử nghiệm." inside the table,
#    reproducing the resulting 3-run split ("Đã t" / "ì" / "m hiểu ...").
# ---------------------------------------------------------------------------
$t = $d.Tables(1)

for ($row = 1; $row -le $t.Rows.Count; $row++) {
    $cell = $t.Cell($row, 3)
    $cellRange = $cell.Range
    $scoped = $d.Range($cellRange.Start, $cellRange.End)
    $found = $scoped.Find.Execute("tiềm")
    if ($found) {
        $s = $scoped.Start

        # "ềm" -> "m"  (drop the "ề", keep a plain "m") as its own run
        $tail = $d.Range($s + 2, $s + 4)
        $tail.Text = "m"
        $tailToggle = $d.Range($s + 2, $s + 3)
        $tailToggle.Font.Bold = $true
        $tailToggle.Font.Bold = $false

        # "i" -> "ì" as its own run
        $mid = $d.Range($s + 1, $s + 2)
        $mid.Text = "ì"
        $mid.Font.Bold = $true
        $mid.Font.Bold = $false
    }
}
